# Applies the cryptos.xlsx price/volume/coin refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-3 character used inside one of the price strings (e.g. 0.0₃0788)
$sub3 = [char]0x2083

$ws.Range("D2").Value = "34.998.71"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.824.42"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "230.88"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.03%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "39.97"
$ws.Range("E8").Value = "  -2.06%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.321"
$ws.Range("E9").Value = "  +5.52%  "
$ws.Range("E10").Value = "  +0.59%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.0993"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "2.089.16"
$ws.Range("E12").Value = "  +0.66%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "11.34"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.826.77"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.668"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "35.010.61"
$ws.Range("E17").Value = "  +0.57%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "69.71"
$ws.Range("E18").Value = "  +1.35%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = [string]::Concat("0.0", $sub3, "0788")
$ws.Range("E19").Value = "  +1.00%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "240.50"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  +3.77%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "4.67"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  +0.11%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "2.27"
$ws.Range("E24").Value = "  +1.90%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "173.67"
$ws.Range("E25").Value = "  +0.79%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "7.83"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("E27").Value = "  +3.86%  "
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "3.99"
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.0551"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +12.13%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "1.82"
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "1.41"
$ws.Range("E36").Value = "  +10.49%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.698"
$ws.Range("E37").Value = "  +4.36%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "92.98"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "1.340.61"
$ws.Range("E39").Value = "  +2.98%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.0194"
$ws.Range("E40").Value = "  +1.64%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "1.00"
$ws.Range("E41").Value = "  +2.96%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "14.67"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "2.43"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "2.26"
$ws.Range("E44").Value = "  -1.76%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "2.75"
$ws.Range("E45").Value = "  +0.41%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "6.27"
$ws.Range("E46").Value = "  +1.57%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.0522"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").Value = "2.006.19"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +0.07%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "0.0671"
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"   # force text so Excel does not coerce this price into a number/date
$cell.Value = "97.98"
$ws.Range("E51").Value = "  -1.48%  "
